$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'311.91"
$ws.Range("E2").Value = "'1.63%"
$ws.Range("D3").Value = "'37.49"
$ws.Range("E3").Value = "'0.97%"
$ws.Range("D4").Value = "'5.124"
$ws.Range("E4").Value = "'1.06%"
$ws.Range("D5").Value = "'0.07873"
$ws.Range("E5").Value = "'1.96%"
$ws.Range("B6").Value = "'GateToken"
$ws.Range("C6").Value = "'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D6").Value = "'4.422"
$ws.Range("E6").Value = "'1.90%"
$ws.Range("B7").Value = "'KuCoinToken"
$ws.Range("C7").Value = "'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D7").Value = "'8.264"
$ws.Range("E7").Value = "'0.92%"
$ws.Range("B8").Value = "'FTXToken"
$ws.Range("C8").Value = "'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D8").Value = "'1.900"
$ws.Range("E8").Value = "'0.64%"
$ws.Range("B9").Value = "'BTSEToken"
$ws.Range("C9").Value = "'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D9").Value = "'2.838"
$ws.Range("E9").Value = "'-10.59%"
$ws.Range("B10").Value = "'MXToken"
$ws.Range("C10").Value = "'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D10").Value = "'0.9203"
$ws.Range("E10").Value = "'0.16%"
$ws.Range("B11").Value = "'LiechtensteinCryptoassetsExchange"
$ws.Range("C11").Value = "'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D11").Value = "'0.1185"
$ws.Range("E11").Value = "'-5.01%"
$ws.Range("B12").Value = "'WazirX"
$ws.Range("C12").Value = "'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D12").Value = "'0.1930"
$ws.Range("E12").Value = "'3.20%"
$ws.Range("B13").Value = "'MandalaExchangeToken"
$ws.Range("C13").Value = "'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D13").Value = "'0.09058"
$ws.Range("E13").Value = "'2.99%"
$ws.Range("B14").Value = "'BitrueCoin"
$ws.Range("C14").Value = "'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D14").Value = "'0.03320"
$ws.Range("E14").Value = "'-1.93%"
$ws.Range("B15").Value = "'BitMartToken"
$ws.Range("C15").Value = "'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D15").Value = "'0.09606"
$ws.Range("E15").Value = "'-0.88%"
$ws.Range("B16").Value = "'BitForexToken"
$ws.Range("C16").Value = "'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D16").Value = "'0.001384"
$ws.Range("E16").Value = "'1.19%"
$ws.Range("B17").Value = "'TigerCash"
$ws.Range("C17").Value = "'https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D17").Value = "'0.005933"
$ws.Range("E17").Value = "'-0.53%"
$ws.Range("B18").Value = "'LEO"
$ws.Range("C18").Value = "'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D18").Value = "'3.554"
$ws.Range("E18").Value = "'-1.00%"
$ws.Range("D19").Value = "'0.3442"
$ws.Range("E19").Value = "'0.97%"
$ws.Range("D20").Value = "'5.259"
$ws.Range("E20").Value = "'4.88%"
$ws.Range("D21").Value = "'0.1284"
$ws.Range("E21").Value = "'1.16%"
$ws.Range("E22").Value = "'3.99%"
$ws.Range("D23").Value = "'0.04350"
$ws.Range("E23").Value = "'0.55%"
$ws.Range("D24").Value = "'0.001247"
$ws.Range("E24").Value = "'3.10%"
$ws.Range("D25").Value = "'0.004661"
$ws.Range("E25").Value = "'10.14%"
$ws.Range("D26").Value = "'0.0001358"
$ws.Range("E26").Value = "'0.63%"
$ws.Range("E27").Value = "'-98.10%"
$ws.Range("D39").Value = "'0.02259"
$ws.Range("E39").Value = "'3.91%"
$ws.Range("D40").Value = "'0.05080"
$ws.Range("E40").Value = "'3.86%"
$ws.Range("D41").Value = "'0.007441"
$ws.Range("E41").Value = "'-3.37%"
$ws.Range("D42").Value = "'0.009050"
$ws.Range("E42").Value = "'-8.49%"
$ws.Range("D43").Value = "'0.1351"
$ws.Range("E43").Value = "'0.93%"
$ws.Range("D44").Value = "'0.001947"
$ws.Range("E44").Value = "'-2.30%"
$ws.Range("D45").Value = "'0.008605"
$ws.Range("E45").Value = "'-12.56%"
$ws.Range("D46").Value = "'0.00006575"
$ws.Range("E46").Value = "'0.73%"
$ws.Range("D47").Value = "'0.00000000749"
$ws.Range("E47").Value = "'-0.11%"
$ws.Range("D48").Value = "'0.003303"
$ws.Range("E48").Value = "'10.13%"
$ws.Range("D49").Value = "'0.0009997"
$ws.Range("E49").Value = "'-23.07%"
$ws.Range("D50").Value = "'0.00002097"
$ws.Range("E50").Value = "'-0.11%"
$ws.Range("D51").Value = "'0.0001997"
$ws.Range("E51").Value = "'-0.11%"
